# "Generate Report for Archive"
#
# The localization status for the b85bd77a-... handoff moved on from
# "Ready for handoff" to "In Translation". That shared string is used on
# all three sheets (Overview!E2/F2 and the per-locale Status column, C2,
# on both the zh-cn and de-de sheets) - update every occurrence so they
# all end up pointing at the same new text.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# "In Translation" is shorter than "Ready for handoff", so the
# auto-fitted Status/zh-cn/de-de columns shrink to match the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
